# Loan RBI, Variable Instalments
# Insert a new (blank) column into the "Repayment schedule" sheet at column N,
# shifting the existing "Late" / "Outstanding" (heading) / "Outstanding" columns
# one position to the right, and make that sheet the active tab/selection.

$wb = $excel.ActiveWorkbook
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column before current column N (14th column), pushing
# N->O, O->P, P->Q for every row (header + the 7 data rows).
$wsRepay.Columns.Item(14).Insert()

# The newly inserted column takes on a width of 11 (same displayed width as
# column M), but without the "best fit" auto-sizing flag.
$wsRepay.Columns.Item(14).ColumnWidth = 10.17

# Switch focus to the Repayment schedule sheet and leave the selection where
# the edit was made.
$wsRepay.Activate() | Out-Null
$wsRepay.Range("N14").Select() | Out-Null
